# "Adjust the 2pc figure"
#
# Slide 4 has a small figure made of vertical dashed "timeline" connectors
# plus numbered callout textboxes (1 / 2 / 3). It is being trimmed from a
# three-timeline figure down to a two-timeline one: the 3rd connector and
# the "3" callout go away, the 2nd connector's slot is taken over by the
# 1st connector (repositioned/lengthened and given a new dash style), and
# the two surviving callouts ("1" and "2") slide over to line up with it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Drop the shapes that fall out of the figure --------------------------
$s.Shapes.Item("Straight Connector 49").Delete()   # 3rd timeline connector
$s.Shapes.Item("Straight Connector 48").Delete()   # 2nd timeline connector (superseded below)
$s.Shapes.Item("TextBox 57").Delete()              # "3" callout

# --- Re-purpose the 1st connector as the new 2nd connector ----------------
$connector = $s.Shapes.Item("Straight Connector 46")
$connector.Name = "Straight Connector 48"
# Target EMU off/ext: x=3059832 y=1366817 cx=0 cy=1383588
# (Left/Top/Height are expressed in points, as PowerPoint's COM model
# requires; the literals below are the point values that round-trip to the
# exact target EMU through that Single-precision conversion.)
$connector.Left = 240.9316535433071
$connector.Top = 107.62338638305665
$connector.Width = 0
$connector.Height = 108.94393700787401
$connector.Line.Weight = 0.75           # -> <a:ln w="9525">
$connector.Line.DashStyle = 6           # msoLineDashDotDot -> prstDash "lgDashDotDot"

# --- Slide the two surviving callouts into their new spots ----------------
$label1 = $s.Shapes.Item("TextBox 55")
$label1.Name = "TextBox 28"
# Target EMU off: x=1691680 y=1366817
$label1.Left = 133.20314960629923
$label1.Top = 107.62338638305665

$label2 = $s.Shapes.Item("TextBox 56")
$label2.Name = "TextBox 30"
# Target EMU off: x=3498106 y=1366817
$label2.Left = 275.4414215087891
$label2.Top = 107.62338638305665
